$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Samples Tab" query, now trimmed: the Tumor / Analyte Type columns
# (smp.sample_tumor_status, smp.sample_type) are dropped from the SELECT list.
$samplesQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs001524' AND sp.gender = 'Male'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

# The workbook's multi-line cell text uses Windows-style CRLF line breaks
# (as produced by Excel's Alt+Enter); normalize our here-string to match.
$samplesQuery = $samplesQuery -replace "`r`n", "`n"
$samplesQuery = $samplesQuery -replace "`n", "`r`n"

# B3 ("SamplesTab" row) now carries the trimmed Samples query. B4
# ("FilesTab" row) already holds the identical Files query text, so it is
# left untouched (re-writing it would only reinsert an identical shared
# string, reordering the shared-strings table unnecessarily).
$ws.Range("B3").Value = $samplesQuery

# Update the view: scroll so row 3 is at the top and select C3 instead of C4.
$ws.Range("C3").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
